$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> column -> new value, derived from the target diff.
$data = @{
    2  = @{ E=3; G=33.47808166666666; H=100.434245; I=0.4880542983452505; J=0.4880542983452505; K=3; M=28.72417333333333; N=86.17251999999999; O=0.4233259107972328; P=0.4233259107972328; Q=961.6302206608219; R=8654.671985947398; S=0.2066060303655075; T=0.2066060303655076 }
    3  = @{ E=3; G=33.47808166666666; H=100.434245; I=0.4880542983452505; J=0.4880542983452505; K=3; M=30.56986233333333; N=91.709587; O=0.4505269713084062; P=0.4505269713084062; Q=1023.42034773409; R=9210.783129606814; S=0.219881624867535; T=0.219881624867535 }
    4  = @{ E=3; G=33.47808166666666; H=100.434245; I=0.4880542983452505; J=0.4880542983452505; K=3; M=8.559531999999999; N=25.678596; O=0.126147117894361; P=0.126147117894361; Q=286.5567113244466; R=2579.01040192002; S=0.06156664311220794; T=0.06156664311220794 }
    5  = @{ E=3; G=23.24776266666667; H=69.74328800000001; I=0.3389134003957588; J=0.3389134003957588; K=3; M=28.72417333333333; N=86.17251999999999; O=0.4233259107972328; P=0.4233259107972328; Q=667.7727644495289; R=6009.95488004576; S=0.1434708239039218; T=0.1434708239039218 }
    6  = @{ E=3; G=23.24776266666667; H=69.74328800000001; I=0.3389134003957588; J=0.3389134003957588; K=3; M=30.56986233333333; N=91.709587; O=0.4505269713084062; P=0.4505269713084062; Q=710.6809042780063; R=6396.128138502057; S=0.1526896278161344; T=0.1526896278161344 }
    7  = @{ E=3; G=23.24776266666667; H=69.74328800000001; I=0.3389134003957588; J=0.3389134003957588; K=3; M=8.559531999999999; N=25.678596; O=0.126147117894361; P=0.126147117894361; Q=198.9899684737387; R=1790.909716263648; S=0.04275294867570256; T=0.04275294867570255 }
    8  = @{ E=3; G=11.86914966666667; H=35.607449; I=0.1730323012589908; J=0.1730323012589908; K=3; M=28.72417333333333; N=86.17251999999999; O=0.4233259107972328; P=0.4233259107972328; Q=340.9315123446088; R=3068.383611101479; S=0.07324905652780346; T=0.07324905652780345 }
    9  = @{ E=3; G=11.86914966666667; H=35.607449; I=0.1730323012589908; J=0.1730323012589908; K=3; M=30.56986233333333; N=91.709587; O=0.4505269713084062; P=0.4505269713084062; Q=362.8382713237291; R=3265.544441913562; S=0.07795571862473687; T=0.07795571862473685 }
    10 = @{ E=3; G=11.86914966666667; H=35.607449; I=0.1730323012589908; J=0.1730323012589908; K=3; M=8.559531999999999; N=25.678596; O=0.126147117894361; P=0.126147117894361; Q=101.5943663846226; R=914.3492974616039; S=0.0218275261064505; T=0.0218275261064505 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
